$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: GaluaPulemet6 -> GaluaPulemet64
$ws.Range("B2").Value = "GaluaPulemet64"
$ws.Range("C2").Value = "GaluaPulemet64@gmail.com"

# Row 3: GaluaPulemet7 -> GaluaPulemet72
$ws.Range("B3").Value = "GaluaPulemet72"
$ws.Range("C3").Value = "GaluaPulemet72@gmail.com"

# Row 4: GaluaPulemet8 -> GaluaPulemet82
$ws.Range("B4").Value = "GaluaPulemet82"
$ws.Range("C4").Value = "GaluaPulemet82@gmail.com"

# Row 5: GaluaPulemet9 -> GaluaPulemet92
$ws.Range("B5").Value = "GaluaPulemet92"
$ws.Range("C5").Value = "GaluaPulemet92@gmail.com"

# Row 6: Magnifikate77@gmail.com -> Magnifikate772@gmail.com
$ws.Range("C6").Value = "Magnifikate772@gmail.com"

# Update selection to I9
$ws.Range("I9").Select()
